# Mise a jour 0.0.5
# - Replace the cascading "+1" formulas in column A (rows 3-15, 17-18) with
#   plain numeric values (keeping the same displayed numbers).
# - Fill in the previously empty A16 cell with its sequence value (15).
# - Because row 16 now participates in the sequence, every following row's
#   number shifts up by one (A17: 15->16, A18: 16->17, A19: 17->18).
# - Move the active selection from E14 to A20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value  = 2
$ws.Range("A4").Value  = 3
$ws.Range("A5").Value  = 4
$ws.Range("A6").Value  = 5
$ws.Range("A7").Value  = 6
$ws.Range("A8").Value  = 7
$ws.Range("A9").Value  = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12
$ws.Range("A14").Value = 13
$ws.Range("A15").Value = 14
$ws.Range("A16").Value = 15
$ws.Range("A17").Value = 16
$ws.Range("A18").Value = 17
$ws.Range("A19").Value = 18

$ws.Range("A20").Select()
